$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column BB: next release vintage (2025-11-25, serial 45986) ---
# Copy header cell formatting (date format with border, from BA1) to BB1
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)
$ws.Range("BB1").Value = 45986

# For every existing quarter (rows 2-82) the newest vintage simply repeats the prior (BA) estimate
$ws.Range("BB2").Value = 0.3540266550812703
$ws.Range("BB3").Value = 0.4017652379482399
$ws.Range("BB4").Value = 0.5660678925619038
$ws.Range("BB5").Value = 0.2814436235376547
$ws.Range("BB6").Value = 0.8315557155661395
$ws.Range("BB7").Value = 1.495825053739068
$ws.Range("BB8").Value = 0.6990997280735201
$ws.Range("BB9").Value = 1.050759988691979
$ws.Range("BB10").Value = 0.3241927986796327
$ws.Range("BB11").Value = 0.3231523476416669
$ws.Range("BB12").Value = 0.8006648126470708
$ws.Range("BB13").Value = 0.1369474440744227
$ws.Range("BB14").Value = 1.366270496737897
$ws.Range("BB15").Value = -0.6985480789094254
$ws.Range("BB16").Value = -0.4148631161428114
$ws.Range("BB17").Value = -2.200693209579313
$ws.Range("BB18").Value = -4.008176425042492
$ws.Range("BB19").Value = 0.3143994836296855
$ws.Range("BB20").Value = 0.8129351329973105
$ws.Range("BB21").Value = 0.7286484790473651
$ws.Range("BB22").Value = 0.662286234074088
$ws.Range("BB23").Value = 2.210734389673945
$ws.Range("BB24").Value = 0.702250570321695
$ws.Range("BB25").Value = 0.6043683783303493
$ws.Range("BB26").Value = 1.531899235856926
$ws.Range("BB27").Value = 0.09102045989541807
$ws.Range("BB28").Value = 0.4274346081797518
$ws.Range("BB29").Value = 0.09961124419930911
$ws.Range("BB30").Value = 0.3110551286058296
$ws.Range("BB31").Value = 0.1249967678526218
$ws.Range("BB32").Value = 0.08642692979871924
$ws.Range("BB33").Value = -0.4125838437329037
$ws.Range("BB34").Value = -0.404671332649869
$ws.Range("BB35").Value = 0.793265767125348
$ws.Range("BB36").Value = 0.2975365817668774
$ws.Range("BB37").Value = 0.4497576285229741
$ws.Range("BB38").Value = 0.715037077548871
$ws.Range("BB39").Value = -0.05661231354093843
$ws.Range("BB40").Value = 0.1888053351092367
$ws.Range("BB41").Value = 0.6124584237519315
$ws.Range("BB42").Value = 0.170162512332567
$ws.Range("BB43").Value = 0.5212263459736306
$ws.Range("BB44").Value = 0.2407365673923465
$ws.Range("BB45").Value = 0.3602440716739608
$ws.Range("BB46").Value = 0.7156188677996056
$ws.Range("BB47").Value = 0.4660496629244335
$ws.Range("BB48").Value = 0.1637213432474738
$ws.Range("BB49").Value = 0.4177287092911968
$ws.Range("BB50").Value = 0.8997744569043959
$ws.Range("BB51").Value = 0.633132069676634
$ws.Range("BB52").Value = 0.735487593389081
$ws.Range("BB53").Value = 0.6069719124519395
$ws.Range("BB54").Value = 0.1320325676681762
$ws.Range("BB55").Value = 0.3861390137996494
$ws.Range("BB56").Value = -0.1032042178152466
$ws.Range("BB57").Value = 0.2066105200339621
$ws.Range("BB58").Value = 0.5
$ws.Range("BB59").Value = -0.2
$ws.Range("BB60").Value = 0.3
$ws.Range("BB61").Value = -0.1
$ws.Range("BB62").Value = -1.995361287679273
$ws.Range("BB63").Value = -9.697717272052344
$ws.Range("BB64").Value = 8.701161067295743
$ws.Range("BB65").Value = 0.5342924547835821
$ws.Range("BB66").Value = -1.684226516424943
$ws.Range("BB67").Value = 2.173700732922356
$ws.Range("BB68").Value = 1.669530332166502
$ws.Range("BB69").Value = -0.3471888372093019
$ws.Range("BB70").Value = 1.023580707979747
$ws.Range("BB71").Value = -0.1208922437305517
$ws.Range("BB72").Value = 0.4746827657805142
$ws.Range("BB73").Value = -0.5372032863913034
$ws.Range("BB74").Value = 0.2711818952007405
$ws.Range("BB75").Value = -0.07397500112315925
$ws.Range("BB76").Value = 0.1388009164387825
$ws.Range("BB77").Value = -0.4989801917293875
$ws.Range("BB78").Value = 0.2365428825421532
$ws.Range("BB79").Value = -0.2955890549112326
$ws.Range("BB80").Value = 0.1051437241507784
$ws.Range("BB81").Value = -0.2005382402049349
$ws.Range("BB82").Value = 0.3062442926496516

# Row 83 (2025 Q1, date 2025-05-15) receives a small revision in the new vintage
$ws.Range("BB83").Value = -0.2099029780610664

# --- New row 84: newest quarter (2025 Q2, date 2025-08-15) first appears in this vintage ---
# Copy date-cell formatting from A83 to A84
$ws.Range("A83").Copy()
$ws.Range("A84").PasteSpecial(-4122)
$ws.Range("A84").Value = 45884
$ws.Range("BB84").Value = 0

$excel.CutCopyMode = 0
